# Apply the data refresh to both "展览" and "全部类型" sheets (identical tables).
$wb = $excel.ActiveWorkbook
$sheetNames = @("展览", "全部类型")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 2: 景德镇·原神X崩铁X崩坏动漫展only
    $ws.Range("B2").NumberFormat = "@"
    $ws.Range("B2").Value = "2024-03-16"
    $ws.Range("B2").Style = "Normal"
    $ws.Range("C2").Value = "景德镇·原神X崩铁X崩坏动漫展only"
    $ws.Range("D2").Value = "陶阳南路188号 晨枫臻品酒店"
    $ws.Range("E2").Value = "2024.03.16 10:00-03.16 17:00"
    $ws.Range("F2").Value = 78
    $ws.Range("G2").Value = 55
    $ws.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=80920"
    $ws.Range("I2").Value = "//i0.hdslb.com/bfs/openplatform/202401/IugBckTp1705469476482.png"

    # Row 3: 江西·ShiningStaR动漫游戏文化节5th
    $ws.Range("B3").NumberFormat = "@"
    $ws.Range("B3").Value = "2024-03-16"
    $ws.Range("B3").Style = "Normal"
    $ws.Range("C3").Value = "江西·ShiningStaR动漫游戏文化节5th"
    $ws.Range("D3").Value = "江西科技学院内 江西科技学院体育馆"
    $ws.Range("E3").Value = "2024.03.16 09:30-03.17 17:00"
    $ws.Range("F3").Value = 3089
    $ws.Range("G3").Value = "不可售"
    $ws.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=81792"
    $ws.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202403/p3TpZeAQ1709544877660.jpeg"

    # Row 4: 上饶·原×铁×崩only
    $ws.Range("B4").NumberFormat = "@"
    $ws.Range("B4").Value = "2024-03-23"
    $ws.Range("B4").Style = "Normal"
    $ws.Range("C4").Value = "上饶·原×铁×崩only"
    $ws.Range("D4").Value = "五三东大道42号 回禾酒店"
    $ws.Range("E4").Value = "2024.03.23 10:00-03.23 17:00"
    $ws.Range("F4").Value = 42
    $ws.Range("G4").Value = 60
    $ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=81103"
    $ws.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202401/pp6c5TsC1705647180602.jpeg"

    # Row 5: 南昌·AP动漫游戏嘉年华
    $ws.Range("B5").NumberFormat = "@"
    $ws.Range("B5").Value = "2024-03-23"
    $ws.Range("B5").Style = "Normal"
    $ws.Range("C5").Value = "南昌·AP动漫游戏嘉年华"
    $ws.Range("D5").Value = "八一桥街道青山南路118号 蓝海会展中心"
    $ws.Range("E5").Value = "2024.03.23 09:00-03.24 17:00"
    $ws.Range("F5").Value = 2537
    $ws.Range("G5").Value = 58.5
    $ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=81232"
    $ws.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202401/NZv97SmS1705912230957.jpeg"

    # Row 6: 南昌·运动番only春季集训（取消）
    $ws.Range("B6").NumberFormat = "@"
    $ws.Range("B6").Value = "2024-03-23"
    $ws.Range("B6").Style = "Normal"
    $ws.Range("C6").Value = "南昌·运动番only春季集训（取消）"
    $ws.Range("D6").Value = "创新三路777号 南昌小飞侠章鱼文化体育公园"
    $ws.Range("E6").Value = "2024.03.23 10:00-03.24 17:00"
    $ws.Range("F6").Value = 186
    $ws.Range("G6").Value = "不可售"
    $ws.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=81950"
    $ws.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202402/bm4uH4qB1708425538357.jpeg"

    # Row 7: 南昌·AP动漫游戏  嘉年华内场票-小N&子音
    $ws.Range("B7").NumberFormat = "@"
    $ws.Range("B7").Value = "2024-03-24"
    $ws.Range("B7").Style = "Normal"
    $ws.Range("C7").Value = "南昌·AP动漫游戏  嘉年华内场票-小N&子音"
    $ws.Range("D7").Value = "八一桥街道青山南路118号 蓝海会展中心"
    $ws.Range("E7").Value = "2024.03.24 09:00-03.24 17:00"
    $ws.Range("F7").Value = 129
    $ws.Range("G7").Value = 218
    $ws.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=81973"
    $ws.Range("I7").Value = "//i0.hdslb.com/bfs/openplatform/202402/zbG5HICL1708504962467.jpeg"

    # Row 8: 鹰潭·宅舞联萌·随舞动漫派对（免费活动)
    $ws.Range("B8").NumberFormat = "@"
    $ws.Range("B8").Value = "2024-03-24"
    $ws.Range("B8").Style = "Normal"
    $ws.Range("C8").Value = "鹰潭·宅舞联萌·随舞动漫派对（免费活动)"
    $ws.Range("D8").Value = "玉清路与象山路交叉口东南角 鹰潭天虹购物中心"
    $ws.Range("E8").Value = "2024.03.24 14:00-03.24 18:00"
    $ws.Range("F8").Value = 4
    $ws.Range("G8").Value = 22.33
    $ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=82434"
    $ws.Range("I8").Value = "//i0.hdslb.com/bfs/openplatform/202403/oj5AEi4W1709713367717.jpeg"

    # Row 9: 南昌·CM01动漫游戏博览会
    $ws.Range("B9").NumberFormat = "@"
    $ws.Range("B9").Value = "2024-03-30"
    $ws.Range("B9").Style = "Normal"
    $ws.Range("C9").Value = "南昌·CM01动漫游戏博览会"
    $ws.Range("D9").Value = "怀玉山大道1315号 南昌绿地国际博览中心"
    $ws.Range("E9").Value = "2024.03.30 10:00-03.31 17:00"
    $ws.Range("F9").Value = 1290
    $ws.Range("G9").Value = 55
    $ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=81691"
    $ws.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202402/9cMJMElF1708938074308.png"

    # Row 10: 鹰潭·原×铁×崩only
    $ws.Range("B10").NumberFormat = "@"
    $ws.Range("B10").Value = "2024-03-30"
    $ws.Range("B10").Style = "Normal"
    $ws.Range("C10").Value = "鹰潭·原×铁×崩only"
    $ws.Range("D10").Value = "南站路24号 回禾酒店(鹰潭火车站店)"
    $ws.Range("E10").Value = "2024.03.30 10:00-03.30 17:00"
    $ws.Range("F10").Value = 36
    $ws.Range("G10").Value = 60
    $ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=81097"
    $ws.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202401/q0AZaXAk1705646244207.jpeg"

    # Row 11: 新余·文旅国漫嘉年华暨BM次元盛典
    $ws.Range("B11").NumberFormat = "@"
    $ws.Range("B11").Value = "2024-03-31"
    $ws.Range("B11").Style = "Normal"
    $ws.Range("C11").Value = "新余·文旅国漫嘉年华暨BM次元盛典"
    $ws.Range("D11").Value = "五一南路与仙女湖大道交叉口西北 老上海风情街白金汉宫"
    $ws.Range("E11").Value = "2024.03.31 10:00-03.31 17:00"
    $ws.Range("F11").Value = 53
    $ws.Range("G11").Value = 60
    $ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=82208"
    $ws.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202403/aXc6vPDP1709547191851.jpeg"

    # Row 12: 景德镇·宅舞联萌·随舞动漫派对（免费活动)
    $ws.Range("B12").NumberFormat = "@"
    $ws.Range("B12").Value = "2024-03-31"
    $ws.Range("B12").Style = "Normal"
    $ws.Range("C12").Value = "景德镇·宅舞联萌·随舞动漫派对（免费活动)"
    $ws.Range("D12").Value = "经二路与纬二路交叉路口 景德镇市宝龙广场"
    $ws.Range("E12").Value = "2024.03.31 14:00-03.31 18:00"
    $ws.Range("F12").Value = 11
    $ws.Range("G12").Value = 22.33
    $ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=82437"
    $ws.Range("I12").Value = "//i0.hdslb.com/bfs/openplatform/202403/zcNNsicc1709714035066.jpeg"

    # Row 13: 南昌·创造力动漫游戏嘉年华1.0
    $ws.Range("B13").NumberFormat = "@"
    $ws.Range("B13").Value = "2024-04-04"
    $ws.Range("B13").Style = "Normal"
    $ws.Range("C13").Value = "南昌·创造力动漫游戏嘉年华1.0"
    $ws.Range("D13").Value = "八一桥街道青山南路118号 蓝海会展中心"
    $ws.Range("E13").Value = "2024.04.04 10:00-04.05 17:00"
    $ws.Range("F13").Value = 1146
    $ws.Range("G13").Value = 39.9
    $ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=82419"
    $ws.Range("I13").Value = "//i2.hdslb.com/bfs/openplatform/202403/dSLjsLyX1709709665435.jpeg"

    # Row 14: 赣州·第三届半夏动漫展
    $ws.Range("B14").NumberFormat = "@"
    $ws.Range("B14").Value = "2024-04-04"
    $ws.Range("B14").Style = "Normal"
    $ws.Range("C14").Value = "赣州·第三届半夏动漫展"
    $ws.Range("D14").Value = "105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心"
    $ws.Range("E14").Value = "2024.04.04 10:00-04.06 17:00"
    $ws.Range("F14").Value = 328
    $ws.Range("G14").Value = 50
    $ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=82235"
    $ws.Range("I14").Value = "//i0.hdslb.com/bfs/openplatform/202403/4DWZWYGm1709278879159.jpeg"

    # Row 15: 赣州·赣次元·归来国风动漫节
    $ws.Range("B15").NumberFormat = "@"
    $ws.Range("B15").Value = "2024-04-04"
    $ws.Range("B15").Style = "Normal"
    $ws.Range("C15").Value = "赣州·赣次元·归来国风动漫节"
    $ws.Range("D15").Value = "客家大道568号文清外国语学校旁 赣州市文清外国语学校国际交流中心"
    $ws.Range("E15").Value = "2024.04.04 10:00-04.04 17:00"
    $ws.Range("F15").Value = 320
    $ws.Range("G15").Value = 40
    $ws.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=82125"
    $ws.Range("I15").Value = "//i1.hdslb.com/bfs/openplatform/202403/fIehikk51709705287036.jpeg"

    # Row 16: 抚州·第七届FZ动漫文化节
    $ws.Range("B16").NumberFormat = "@"
    $ws.Range("B16").Value = "2024-04-05"
    $ws.Range("B16").Style = "Normal"
    $ws.Range("C16").Value = "抚州·第七届FZ动漫文化节"
    $ws.Range("D16").Value = "迎宾大道288号 凤凰世纪名都大酒店"
    $ws.Range("E16").Value = "2024.04.05 09:30-04.05 17:00"
    $ws.Range("F16").Value = 27
    $ws.Range("G16").Value = 50
    $ws.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=82381"
    $ws.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202403/Y725SN0G1709694367526.jpeg"

    # Row 17: 萍乡·2024DDS国漫盛典
    $ws.Range("B17").NumberFormat = "@"
    $ws.Range("B17").Value = "2024-04-06"
    $ws.Range("B17").Style = "Normal"
    $ws.Range("C17").Value = "萍乡·2024DDS国漫盛典"
    $ws.Range("D17").Value = "凤凰街迎宾路18号 鸿凯大酒店"
    $ws.Range("E17").Value = "2024.04.06 10:00-04.06 17:00"
    $ws.Range("F17").Value = 28
    $ws.Range("G17").Value = 30
    $ws.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=82413"
    $ws.Range("I17").Value = "//i1.hdslb.com/bfs/openplatform/202403/Rbu1xyFr1709707847098.jpeg"

    # Row 18: 南昌·原X穹X崩only
    $ws.Range("B18").NumberFormat = "@"
    $ws.Range("B18").Value = "2024-04-13"
    $ws.Range("B18").Style = "Normal"
    $ws.Range("C18").Value = "南昌·原X穹X崩only"
    $ws.Range("D18").Value = "丰和北大道299号 新吉花园酒店"
    $ws.Range("E18").Value = "2024.04.13 10:00-04.13 17:00"
    $ws.Range("F18").Value = 104
    $ws.Range("G18").Value = 65
    $ws.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=80807"
    $ws.Range("I18").Value = "//i0.hdslb.com/bfs/openplatform/202402/kfK13XvH1709202705153.jpeg"

    # Row 19: 南昌·第二届漫拥动漫嘉年华mini
    $ws.Range("B19").NumberFormat = "@"
    $ws.Range("B19").Value = "2024-04-13"
    $ws.Range("B19").Style = "Normal"
    $ws.Range("C19").Value = "南昌·第二届漫拥动漫嘉年华mini"
    $ws.Range("D19").Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
    $ws.Range("E19").Value = "2024.04.13 10:00-04.14 18:00"
    $ws.Range("F19").Value = 68
    $ws.Range("G19").Value = 39.9
    $ws.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=82210"
    $ws.Range("I19").Value = "//i0.hdslb.com/bfs/openplatform/202402/KYd0bfk11709203777701.png"

    # Row 20: 南昌·DSL国风动漫游戏嘉年华
    $ws.Range("B20").NumberFormat = "@"
    $ws.Range("B20").Value = "2024-04-20"
    $ws.Range("B20").Style = "Normal"
    $ws.Range("C20").Value = "南昌·DSL国风动漫游戏嘉年华"
    $ws.Range("D20").Value = "沿江北路69号 瑞颐大酒店"
    $ws.Range("E20").Value = "2024.04.20 09:00-04.21 17:00"
    $ws.Range("F20").Value = 84
    $ws.Range("G20").Value = 35
    $ws.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=82107"
    $ws.Range("I20").Value = "//i0.hdslb.com/bfs/openplatform/202402/QDlumVb41708943318282.jpeg"

    # Row 21: 南昌·New World国潮动漫博览会
    $ws.Range("B21").NumberFormat = "@"
    $ws.Range("B21").Value = "2024-04-20"
    $ws.Range("B21").Style = "Normal"
    $ws.Range("C21").Value = "南昌·New World国潮动漫博览会"
    $ws.Range("D21").Value = "怀玉山大道1315号 南昌绿地国际博览中心"
    $ws.Range("E21").Value = "2024.04.20 09:30-04.21 17:00"
    $ws.Range("F21").Value = 2296
    $ws.Range("G21").Value = 60
    $ws.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=82411"
    $ws.Range("I21").Value = "//i1.hdslb.com/bfs/openplatform/202403/xbYbLXc81709707724935.jpeg"

    # Row 22: 南昌·代号鸢盛花行only
    $ws.Range("B22").NumberFormat = "@"
    $ws.Range("B22").Value = "2024-04-20"
    $ws.Range("B22").Style = "Normal"
    $ws.Range("C22").Value = "南昌·代号鸢盛花行only"
    $ws.Range("D22").Value = "民德路411号 东方豪景花园酒店(民德路店)"
    $ws.Range("E22").Value = "2024.04.20 09:30-04.20 17:30"
    $ws.Range("F22").Value = 19
    $ws.Range("G22").Value = 78
    $ws.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=82529"
    $ws.Range("I22").Value = "//i1.hdslb.com/bfs/openplatform/202403/TJ8iC73c1709804909450.png"

    # Row 23: 九江·第三届ACD动漫游戏嘉年华
    $ws.Range("B23").NumberFormat = "@"
    $ws.Range("B23").Value = "2024-05-01"
    $ws.Range("B23").Style = "Normal"
    $ws.Range("C23").Value = "九江·第三届ACD动漫游戏嘉年华"
    $ws.Range("D23").Value = "九瑞大道与重庆路交汇处西南角 九江国际会展中心"
    $ws.Range("E23").Value = "2024.05.01 09:00-05.02 17:00"
    $ws.Range("F23").Value = 267
    $ws.Range("G23").Value = 39.9
    $ws.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=82464"
    $ws.Range("I23").Value = "//i0.hdslb.com/bfs/openplatform/202403/HjMMyP3a1709780146797.jpeg"

    # Drop the now-stale last row (23 events remain instead of 24) and let the used range shrink.
    $ws.Rows.Item(24).Delete()
}